# Add date-only support: split the existing "Date Time" column's
# formatting from a brand-new "Date" (date-only) column.
#
# Before:  A=Number B=String C=Date Time D=Enum  E=Bool
# After:   A=Number B=String C=Date Time D=Date  E=Enum  F=Bool
#
# Column D becomes the new "Date" column (date-only values); the old
# D/E (Enum/Bool) header text and the row-3 sample values shift one
# column to the right (into E/F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 data values ----------------------------------------------------
# C3: stays the "Date Time" column, but now carries a real time-of-day
# component, with an explicit date+time display format.
$ws.Range("C3").Value = 43831.4271412037
$ws.Range("C3").NumberFormat = "yyyy-MM-dd HH:mm:ss"

# D3: brand-new "Date" column - date-only value/format.
$ws.Range("D3").Value = 43831
$ws.Range("D3").NumberFormat = "yyyy-MM-dd"

# --- Row 1 headers ---------------------------------------------------------
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Enum"
$ws.Range("F1").Value = "Bool"

# --- Remaining row 3 values (old D3/E3 shifted right into E3/F3) ----------
$ws.Range("E3").Value = "Value"
$ws.Range("F3").Value = $true
$ws.Range("F3").HorizontalAlignment = -4131  # xlLeft
$ws.Range("F3").VerticalAlignment = -4160    # xlTop
$ws.Range("F3").WrapText = $true

# --- Row 2 (blank data row) - give the new F2 cell the same blank/styled
# look as the rest of the row.
$ws.Range("F2").HorizontalAlignment = -4131  # xlLeft
$ws.Range("F2").VerticalAlignment = -4160    # xlTop
$ws.Range("F2").WrapText = $true

# --- Column widths: D keeps the new Date column's own width, E/F take on
# the widths that used to belong to D/E (Enum/Bool), shifted right.
$ws.Range("D1").ColumnWidth = 8.5
$ws.Range("E1").ColumnWidth = 9.5
$ws.Range("F1").ColumnWidth = 8.333333333333332

# --- AutoFilter / used range now spans through column F -------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:F3").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name created by AutoFilter needs to
# track the new range too.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$3"
    }
}
